# Auto-generated edit script applying numeric updates per the diff
# (Famfrit_Profits workbook -> market-data columns H..N across ALC/ARM/BSM/CRP/CUL/GSM/LTW)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 638.7778
$ws.Range("I18").Value = 638.7778
$ws.Range("K18").Value = 638.7778
$ws.Range("M18").Value = -354.7778

$ws.Range("H51").Value = 3977.3572
$ws.Range("I51").Value = 2765
$ws.Range("J51").Value = 4551.6313
$ws.Range("K51").Value = 2765
$ws.Range("L51").Value = 4551.6313
$ws.Range("M51").Value = -2281
$ws.Range("N51").Value = -5519.6313

$ws.Range("H55").Value = 176
$ws.Range("J55").Value = 265.33334
$ws.Range("L55").Value = 265.33334
$ws.Range("N55").Value = -693.33334

$ws.Range("H64").Value = 7363
$ws.Range("I64").Value = 7099.3
$ws.Range("K64").Value = 7099.3
$ws.Range("M64").Value = -6851.3

$ws.Range("H67").Value = 7363
$ws.Range("I67").Value = 7099.3
$ws.Range("K67").Value = 7099.3
$ws.Range("M67").Value = -6241.3

$ws.Range("H116").Value = 5175.9414
$ws.Range("I116").Value = 4566.4165
$ws.Range("J116").Value = 6638.8
$ws.Range("K116").Value = 4566.4165
$ws.Range("L116").Value = 6638.8
$ws.Range("M116").Value = -1124.4165
$ws.Range("N116").Value = -13522.8

$ws.Range("H132").Value = 3497.9355
$ws.Range("I132").Value = 3439.1724
$ws.Range("K132").Value = 10317.5172
$ws.Range("M132").Value = -7787.5172

$ws.Range("H137").Value = 6604.625
$ws.Range("I137").Value = 4651.3076
$ws.Range("J137").Value = 8913.091
$ws.Range("K137").Value = 13953.9228
$ws.Range("L137").Value = 26739.273
$ws.Range("M137").Value = -11403.9228
$ws.Range("N137").Value = -31839.273

$ws.Range("H138").Value = 38463316
$ws.Range("I138").Value = 1230.5714
$ws.Range("J138").Value = 83335750
$ws.Range("K138").Value = 3691.7142
$ws.Range("L138").Value = 250007250
$ws.Range("M138").Value = 1448.2858
$ws.Range("N138").Value = -250017530

$ws.Range("H141").Value = 1104.5143
$ws.Range("I141").Value = 983.2727
$ws.Range("K141").Value = 2949.8181
$ws.Range("M141").Value = 2230.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10420029
$ws.Range("I32").Value = 11907760
$ws.Range("K32").Value = 11907760
$ws.Range("M32").Value = -11907473

$ws.Range("H61").Value = 32262118
$ws.Range("I61").Value = 41669550
$ws.Range("K61").Value = 41669550
$ws.Range("M61").Value = -41669338

$ws.Range("H136").Value = 32262118
$ws.Range("I136").Value = 41669550
$ws.Range("K136").Value = 125008650
$ws.Range("M136").Value = -125006100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1802.25
$ws.Range("J64").Value = 1841
$ws.Range("L64").Value = 1841
$ws.Range("N64").Value = -2291

$ws.Range("H67").Value = 1802.25
$ws.Range("J67").Value = 1841
$ws.Range("L67").Value = 1841
$ws.Range("N67").Value = -3401

$ws.Range("H75").Value = 83999.8
$ws.Range("I75").Value = 14999.5
$ws.Range("K75").Value = 14999.5
$ws.Range("M75").Value = -14063.5

$ws.Range("H78").Value = 83999.8
$ws.Range("I78").Value = 14999.5
$ws.Range("K78").Value = 44998.5
$ws.Range("M78").Value = -40318.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30309482
$ws.Range("I31").Value = 3858.8667
$ws.Range("J31").Value = 55564170
$ws.Range("K31").Value = 3858.8667
$ws.Range("L31").Value = 55564170
$ws.Range("M31").Value = -3563.8667
$ws.Range("N31").Value = -55564760

$ws.Range("H34").Value = 30309482
$ws.Range("I34").Value = 3858.8667
$ws.Range("J34").Value = 55564170
$ws.Range("K34").Value = 3858.8667
$ws.Range("L34").Value = 55564170
$ws.Range("M34").Value = -3656.8667
$ws.Range("N34").Value = -55564574

$ws.Range("H62").Value = 2769.5715
$ws.Range("I62").Value = 2831.1667
$ws.Range("K62").Value = 2831.1667
$ws.Range("M62").Value = -2207.1667

$ws.Range("H65").Value = 2769.5715
$ws.Range("I65").Value = 2831.1667
$ws.Range("K65").Value = 14155.8335
$ws.Range("M65").Value = -11035.8335

$ws.Range("H134").Value = 4680.8057
$ws.Range("I134").Value = 4484.654
$ws.Range("J134").Value = 5190.8
$ws.Range("K134").Value = 13453.962
$ws.Range("L134").Value = 15572.4
$ws.Range("M134").Value = -10918.962
$ws.Range("N134").Value = -20642.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 52.25
$ws.Range("I2").Value = 60.63158
$ws.Range("J2").Value = 34.555557
$ws.Range("K2").Value = 363.78948
$ws.Range("L2").Value = 207.333342
$ws.Range("M2").Value = -250.78948
$ws.Range("N2").Value = -433.333342

$ws.Range("H12").Value = 573.5714
$ws.Range("I12").Value = 340.33334
$ws.Range("J12").Value = 637.1818
$ws.Range("K12").Value = 1021.00002
$ws.Range("L12").Value = 1911.5454
$ws.Range("M12").Value = -848.0000200000001
$ws.Range("N12").Value = -2257.5454

$ws.Range("H34").Value = 2259.5557
$ws.Range("I34").Value = 126.28571
$ws.Range("J34").Value = 3617.0908
$ws.Range("K34").Value = 378.85713
$ws.Range("L34").Value = 10851.2724
$ws.Range("M34").Value = -294.85713
$ws.Range("N34").Value = -11019.2724

$ws.Range("H68").Value = 977.6
$ws.Range("J68").Value = 998.25
$ws.Range("L68").Value = 2994.75
$ws.Range("N68").Value = -4616.75

$ws.Range("H71").Value = 977.6
$ws.Range("J71").Value = 998.25
$ws.Range("L71").Value = 8984.25
$ws.Range("N71").Value = -17096.25

$ws.Range("H107").Value = 781.88
$ws.Range("J107").Value = 848.6667
$ws.Range("L107").Value = 2546.0001
$ws.Range("N107").Value = -6386.0001

$ws.Range("H131").Value = 1636.2222
$ws.Range("J131").Value = 1715.6154
$ws.Range("L131").Value = 5146.8462
$ws.Range("N131").Value = -15226.8462

$ws.Range("H132").Value = 3032958
$ws.Range("I132").Value = 1605.3334
$ws.Range("K132").Value = 14448.0006
$ws.Range("M132").Value = -11918.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8250
$ws.Range("I5").Value = 7250
$ws.Range("J5").Value = 8750
$ws.Range("K5").Value = 7250
$ws.Range("L5").Value = 8750
$ws.Range("M5").Value = -7138
$ws.Range("N5").Value = -8974

$ws.Range("H132").Value = 1450.3793
$ws.Range("I132").Value = 1381.96
$ws.Range("K132").Value = 4145.88
$ws.Range("M132").Value = -1615.88

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4485.923
$ws.Range("I7").Value = 3878.6667
$ws.Range("J7").Value = 5852.25
$ws.Range("K7").Value = 3878.6667
$ws.Range("L7").Value = 5852.25
$ws.Range("M7").Value = -3766.6667
$ws.Range("N7").Value = -6076.25

$ws.Range("H16").Value = 3161.85
$ws.Range("I16").Value = 3179.8333
$ws.Range("K16").Value = 3179.8333
$ws.Range("M16").Value = -3009.8333

$ws.Range("H43").Value = 24000
$ws.Range("I43").Value = 24000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 24000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -23807
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 1835.2354
$ws.Range("I46").Value = 799.9286
$ws.Range("K46").Value = 799.9286
$ws.Range("M46").Value = -611.9286

$ws.Range("H99").Value = 39475.223
$ws.Range("I99").Value = 33198.6
$ws.Range("J99").Value = 47321
$ws.Range("K99").Value = 33198.6
$ws.Range("L99").Value = 47321
$ws.Range("M99").Value = -30203.6
$ws.Range("N99").Value = -53311

$ws.Range("H126").Value = 4485.923
$ws.Range("I126").Value = 3878.6667
$ws.Range("J126").Value = 5852.25
$ws.Range("K126").Value = 11636.0001
$ws.Range("L126").Value = 17556.75
$ws.Range("M126").Value = -9166.000100000001
$ws.Range("N126").Value = -22496.75
